$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 50) with date values, continuing the existing table.
$ws.Range("A50").Value = 43348
$ws.Range("B50").Value = 43352

# Match formatting (date number format + borders) of the row above.
$ws.Range("A49:B49").Copy()
$ws.Range("A50:B50").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Scroll the view and select the new cell, matching the target sheetView state.
$ws.Range("B50").Select()
$excel.ActiveWindow.ScrollRow = 17
